$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '66.497.02'
$c.Style = 'Normal'

$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  -1.20%  '
$c.Style = 'Normal'

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.584.99'
$c.Style = 'Normal'

$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -1.90%  '
$c.Style = 'Normal'

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '583.62'
$c.Style = 'Normal'

$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -1.56%  '
$c.Style = 'Normal'

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '166.24'
$c.Style = 'Normal'

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -0.63%  '
$c.Style = 'Normal'

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.Style = 'Normal'

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.527'
$c.Style = 'Normal'

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  -1.65%  '
$c.Style = 'Normal'

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.584.82'
$c.Style = 'Normal'

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  -1.86%  '
$c.Style = 'Normal'

$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -4.11%  '
$c.Style = 'Normal'

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +0.30%  '
$c.Style = 'Normal'

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.356'
$c.Style = 'Normal'

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -1.85%  '
$c.Style = 'Normal'

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '5.17'
$c.Style = 'Normal'

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  -1.38%  '
$c.Style = 'Normal'

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '26.74'
$c.Style = 'Normal'

$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -3.76%  '
$c.Style = 'Normal'

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '3.054.44'
$c.Style = 'Normal'

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -1.90%  '
$c.Style = 'Normal'

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.0000178'
$c.Style = 'Normal'

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -2.52%  '
$c.Style = 'Normal'

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '66.490.64'
$c.Style = 'Normal'

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -0.97%  '
$c.Style = 'Normal'

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.541.96'
$c.Style = 'Normal'

$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -3.09%  '
$c.Style = 'Normal'

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '11.43'
$c.Style = 'Normal'

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -6.37%  '
$c.Style = 'Normal'

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.72'
$c.Style = 'Normal'

$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -4.87%  '
$c.Style = 'Normal'

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '352.29'
$c.Style = 'Normal'

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -2.34%  '
$c.Style = 'Normal'

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.23'
$c.Style = 'Normal'

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -3.03%  '
$c.Style = 'Normal'

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '4.61'
$c.Style = 'Normal'

$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -2.09%  '
$c.Style = 'Normal'

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +0.08%  '
$c.Style = 'Normal'

$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  -3.88%  '
$c.Style = 'Normal'

$c = $ws.Range('B26')
$c.NumberFormat = '@'
$c.Value = 'Litecoin'
$c.Style = 'Normal'

$c = $ws.Range('C26')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c.Style = 'Normal'

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '68.71'
$c.Style = 'Normal'

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -2.67%  '
$c.Style = 'Normal'

$c = $ws.Range('B27')
$c.NumberFormat = '@'
$c.Value = 'Aptos'
$c.Style = 'Normal'

$c = $ws.Range('C27')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c.Style = 'Normal'

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.99'
$c.Style = 'Normal'

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -8.44%  '
$c.Style = 'Normal'

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.715.48'
$c.Style = 'Normal'

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -1.84%  '
$c.Style = 'Normal'

$c = $ws.Range('B29')
$c.NumberFormat = '@'
$c.Value = 'PEPE'
$c.Style = 'Normal'

$c = $ws.Range('C29')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c.Style = 'Normal'

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.0₃0986'
$c.Style = 'Normal'

$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -2.92%  '
$c.Style = 'Normal'

$c = $ws.Range('B30')
$c.NumberFormat = '@'
$c.Value = 'Bittensor'
$c.Style = 'Normal'

$c = $ws.Range('C30')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c.Style = 'Normal'

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '533.95'
$c.Style = 'Normal'

$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -4.16%  '
$c.Style = 'Normal'

$c = $ws.Range('B31')
$c.NumberFormat = '@'
$c.Value = 'InternetComputer(DFINITY)'
$c.Style = 'Normal'

$c = $ws.Range('C31')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c.Style = 'Normal'

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '8.01'
$c.Style = 'Normal'

$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +0.55%  '
$c.Style = 'Normal'

$c = $ws.Range('B32')
$c.NumberFormat = '@'
$c.Value = 'Fetch.AI'
$c.Style = 'Normal'

$c = $ws.Range('C32')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c.Style = 'Normal'

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.33'
$c.Style = 'Normal'

$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -3.57%  '
$c.Style = 'Normal'

$c = $ws.Range('B33')
$c.NumberFormat = '@'
$c.Value = 'PancakeSwap'
$c.Style = 'Normal'

$c = $ws.Range('C33')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c.Style = 'Normal'

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.85'
$c.Style = 'Normal'

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -3.02%  '
$c.Style = 'Normal'

$c = $ws.Range('B34')
$c.NumberFormat = '@'
$c.Value = 'Kaspa'
$c.Style = 'Normal'

$c = $ws.Range('C34')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c.Style = 'Normal'

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.132'
$c.Style = 'Normal'

$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  -2.36%  '
$c.Style = 'Normal'

$c = $ws.Range('B35')
$c.NumberFormat = '@'
$c.Value = 'FirstDigitalUSD'
$c.Style = 'Normal'

$c = $ws.Range('C35')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c.Style = 'Normal'

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'

$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +0.12%  '
$c.Style = 'Normal'

$c = $ws.Range('B36')
$c.NumberFormat = '@'
$c.Value = 'ImmutableX'
$c.Style = 'Normal'

$c = $ws.Range('C36')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c.Style = 'Normal'

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.46'
$c.Style = 'Normal'

$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -3.66%  '
$c.Style = 'Normal'

$c = $ws.Range('B37')
$c.NumberFormat = '@'
$c.Value = 'Monero'
$c.Style = 'Normal'

$c = $ws.Range('C37')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c.Style = 'Normal'

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '156.91'
$c.Style = 'Normal'

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -0.50%  '
$c.Style = 'Normal'

$c = $ws.Range('B38')
$c.NumberFormat = '@'
$c.Value = 'EthereumClassic'
$c.Style = 'Normal'

$c = $ws.Range('C38')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c.Style = 'Normal'

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '18.75'
$c.Style = 'Normal'

$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -2.44%  '
$c.Style = 'Normal'

$c = $ws.Range('B39')
$c.NumberFormat = '@'
$c.Value = 'PolygonEcosystemToken'
$c.Style = 'Normal'

$c = $ws.Range('C39')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c.Style = 'Normal'

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.361'
$c.Style = 'Normal'

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -1.82%  '
$c.Style = 'Normal'

$c = $ws.Range('B40')
$c.NumberFormat = '@'
$c.Value = 'WhiteBITCoin'
$c.Style = 'Normal'

$c = $ws.Range('C40')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c.Style = 'Normal'

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '18.25'
$c.Style = 'Normal'

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +1.80%  '
$c.Style = 'Normal'

$c = $ws.Range('B41')
$c.NumberFormat = '@'
$c.Value = 'Stacks'
$c.Style = 'Normal'

$c = $ws.Range('C41')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c.Style = 'Normal'

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '1.79'
$c.Style = 'Normal'

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -0.89%  '
$c.Style = 'Normal'

$c = $ws.Range('B42')
$c.NumberFormat = '@'
$c.Value = 'RenderToken'
$c.Style = 'Normal'

$c = $ws.Range('C42')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$c.Style = 'Normal'

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '5.11'
$c.Style = 'Normal'

$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -2.22%  '
$c.Style = 'Normal'

$c = $ws.Range('B43')
$c.NumberFormat = '@'
$c.Value = 'USDe'
$c.Style = 'Normal'

$c = $ws.Range('C43')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c.Style = 'Normal'

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.Style = 'Normal'

$c = $ws.Range('B44')
$c.NumberFormat = '@'
$c.Value = 'dogwifhat'
$c.Style = 'Normal'

$c = $ws.Range('C44')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c.Style = 'Normal'

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.42'
$c.Style = 'Normal'

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -2.65%  '
$c.Style = 'Normal'

$c = $ws.Range('B45')
$c.NumberFormat = '@'
$c.Value = 'BabyDogeCoin'
$c.Style = 'Normal'

$c = $ws.Range('C45')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c.Style = 'Normal'

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0₆0286'
$c.Style = 'Normal'

$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -4.73%  '
$c.Style = 'Normal'

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '149.37'
$c.Style = 'Normal'

$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -1.91%  '
$c.Style = 'Normal'

$c = $ws.Range('B47')
$c.NumberFormat = '@'
$c.Value = 'ARBITRUM'
$c.Style = 'Normal'

$c = $ws.Range('C47')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c.Style = 'Normal'

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.567'
$c.Style = 'Normal'

$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -3.66%  '
$c.Style = 'Normal'

$c = $ws.Range('B48')
$c.NumberFormat = '@'
$c.Value = 'Filecoin'
$c.Style = 'Normal'

$c = $ws.Range('C48')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c.Style = 'Normal'

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '3.72'
$c.Style = 'Normal'

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -2.66%  '
$c.Style = 'Normal'

$c = $ws.Range('B49')
$c.NumberFormat = '@'
$c.Value = 'Optimism'
$c.Style = 'Normal'

$c = $ws.Range('C49')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$c.Style = 'Normal'

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.71'
$c.Style = 'Normal'

$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -1.51%  '
$c.Style = 'Normal'

$c = $ws.Range('B50')
$c.NumberFormat = '@'
$c.Value = 'Cronos'
$c.Style = 'Normal'

$c = $ws.Range('C50')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c.Style = 'Normal'

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.0759'
$c.Style = 'Normal'

$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -1.64%  '
$c.Style = 'Normal'

$c = $ws.Range('B51')
$c.NumberFormat = '@'
$c.Value = 'Mantle'
$c.Style = 'Normal'

$c = $ws.Range('C51')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c.Style = 'Normal'

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.597'
$c.Style = 'Normal'

$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -1.40%  '
$c.Style = 'Normal'

